$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 800
$ws.Range("I13").Value = 800
$ws.Range("K13").Value = 800
$ws.Range("M13").Value = -631
$ws.Range("H21").Value = 838.1429000000001
$ws.Range("I21").Value = 838.1429000000001
$ws.Range("K21").Value = 838.1429000000001
$ws.Range("M21").Value = -370.1429000000001
$ws.Range("H23").Value = 838.1429000000001
$ws.Range("I23").Value = 838.1429000000001
$ws.Range("K23").Value = 838.1429000000001
$ws.Range("M23").Value = -604.1429000000001
$ws.Range("H41").Value = 199.25
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H44").Value = 4799.467
$ws.Range("I44").Value = 4571.4287
$ws.Range("J44").Value = 4999
$ws.Range("K44").Value = 4571.4287
$ws.Range("L44").Value = 4999
$ws.Range("M44").Value = -4109.4287
$ws.Range("N44").Value = -5923
$ws.Range("H74").Value = 4496.5
$ws.Range("I74").Value = 4496.5
$ws.Range("K74").Value = 4496.5
$ws.Range("M74").Value = -3560.5
$ws.Range("H77").Value = 4496.5
$ws.Range("I77").Value = 4496.5
$ws.Range("K77").Value = 22482.5
$ws.Range("M77").Value = -17802.5
$ws.Range("H127").Value = 7713.8335
$ws.Range("I127").Value = 8357.799999999999
$ws.Range("J127").Value = 4494
$ws.Range("K127").Value = 25073.4
$ws.Range("L127").Value = 13482
$ws.Range("M127").Value = -20113.4
$ws.Range("N127").Value = -23402
$ws.Range("H132").Value = 1430.1177
$ws.Range("I132").Value = 1430.1177
$ws.Range("K132").Value = 4290.3531
$ws.Range("M132").Value = -1760.3531
$ws.Range("H137").Value = 1524.125
$ws.Range("I137").Value = 1047.5
$ws.Range("K137").Value = 3142.5
$ws.Range("M137").Value = -592.5
$ws.Range("H138").Value = 3480
$ws.Range("I138").Value = 3182.4644
$ws.Range("J138").Value = 7645.5
$ws.Range("K138").Value = 9547.393199999999
$ws.Range("L138").Value = 22936.5
$ws.Range("M138").Value = -4407.393199999999
$ws.Range("N138").Value = -33216.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7312.591
$ws.Range("I32").Value = 5193.5674
$ws.Range("K32").Value = 5193.5674
$ws.Range("M32").Value = -4906.5674
$ws.Range("H63").Value = 8125
$ws.Range("I63").Value = 2750
$ws.Range("K63").Value = 2750
$ws.Range("M63").Value = -2064
$ws.Range("H66").Value = 8125
$ws.Range("I66").Value = 2750
$ws.Range("K66").Value = 13750
$ws.Range("M66").Value = -10318
$ws.Range("H74").Value = 700.75
$ws.Range("J74").Value = 494.66666
$ws.Range("L74").Value = 494.66666
$ws.Range("N74").Value = -2242.66666
$ws.Range("H77").Value = 700.75
$ws.Range("J77").Value = 494.66666
$ws.Range("L77").Value = 2473.3333
$ws.Range("N77").Value = -11209.3333
$ws.Range("H132").Value = 41649.6
$ws.Range("I132").Value = 2749.5
$ws.Range("K132").Value = 8248.5
$ws.Range("M132").Value = -5718.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2811.5
$ws.Range("I105").Value = 2530.25
$ws.Range("K105").Value = 2530.25
$ws.Range("M105").Value = -783.25
$ws.Range("H107").Value = 2385.6667
$ws.Range("I107").Value = 2385.6667
$ws.Range("K107").Value = 2385.6667
$ws.Range("M107").Value = -465.6667000000002

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 132.94737
$ws.Range("I7").Value = 61.714287
$ws.Range("J7").Value = 332.4
$ws.Range("K7").Value = 61.714287
$ws.Range("L7").Value = 332.4
$ws.Range("M7").Value = 51.285713
$ws.Range("N7").Value = -558.4
$ws.Range("H31").Value = 3856.4285
$ws.Range("I31").Value = 2199
$ws.Range("K31").Value = 2199
$ws.Range("M31").Value = -1904
$ws.Range("H34").Value = 3856.4285
$ws.Range("I34").Value = 2199
$ws.Range("K34").Value = 2199
$ws.Range("M34").Value = -1997
$ws.Range("H58").Value = 2428.3704
$ws.Range("I58").Value = 1215.6842
$ws.Range("K58").Value = 1215.6842
$ws.Range("M58").Value = -1012.6842
$ws.Range("H94").Value = 612.75
$ws.Range("I94").Value = 655.75
$ws.Range("J94").Value = 569.75
$ws.Range("K94").Value = 655.75
$ws.Range("L94").Value = 569.75
$ws.Range("M94").Value = -204.75
$ws.Range("N94").Value = -1471.75
$ws.Range("H132").Value = 2462
$ws.Range("I132").Value = 2230.75
$ws.Range("K132").Value = 6692.25
$ws.Range("M132").Value = -4162.25
$ws.Range("H134").Value = 2706.8928
$ws.Range("I134").Value = 2569.842
$ws.Range("J134").Value = 2996.2222
$ws.Range("K134").Value = 7709.526
$ws.Range("L134").Value = 8988.6666
$ws.Range("M134").Value = -5174.526
$ws.Range("N134").Value = -14058.6666
$ws.Range("H136").Value = 2428.3704
$ws.Range("I136").Value = 1215.6842
$ws.Range("K136").Value = 3647.0526
$ws.Range("M136").Value = -1097.0526

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 574.75
$ws.Range("I6").Value = 666.3333
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 1998.9999
$ws.Range("L6").Value = 900
$ws.Range("M6").Value = -1885.9999
$ws.Range("N6").Value = -1126
$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 500
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = ""
$ws.Range("N36").Value = -1331
$ws.Range("H38").Value = 6795.4
$ws.Range("J38").Value = 71.166664
$ws.Range("L38").Value = 213.499992
$ws.Range("N38").Value = -907.499992
$ws.Range("H55").Value = 4009.5
$ws.Range("J55").Value = 4009.5
$ws.Range("L55").Value = 12028.5
$ws.Range("N55").Value = -12382.5
$ws.Range("H92").Value = 296.92307
$ws.Range("I92").Value = 290.7143
$ws.Range("K92").Value = 872.1428999999999
$ws.Range("M92").Value = 375.8571000000001
$ws.Range("H121").Value = 397.33334
$ws.Range("J121").Value = 395.5
$ws.Range("L121").Value = 1186.5
$ws.Range("N121").Value = -3806.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4230.5
$ws.Range("I80").Value = 2991.8572
$ws.Range("J80").Value = 5193.8887
$ws.Range("K80").Value = 2991.8572
$ws.Range("L80").Value = 5193.8887
$ws.Range("M80").Value = -1993.8572
$ws.Range("N80").Value = -7189.8887
$ws.Range("H83").Value = 4230.5
$ws.Range("I83").Value = 2991.8572
$ws.Range("J83").Value = 5193.8887
$ws.Range("K83").Value = 14959.286
$ws.Range("L83").Value = 25969.4435
$ws.Range("M83").Value = -9967.286
$ws.Range("N83").Value = -35953.4435
$ws.Range("H102").Value = 3677
$ws.Range("I102").Value = 2023.25
$ws.Range("K102").Value = 2023.25
$ws.Range("M102").Value = -401.25
$ws.Range("H123").Value = 34999.332
$ws.Range("J123").Value = 34999.332
$ws.Range("L123").Value = 34999.332
$ws.Range("N123").Value = -39899.332

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 32250
$ws.Range("I75").Value = 32250
$ws.Range("K75").Value = 32250
$ws.Range("M75").Value = -31314
$ws.Range("H78").Value = 32250
$ws.Range("I78").Value = 32250
$ws.Range("K78").Value = 96750
$ws.Range("M78").Value = -92070
$ws.Range("H82").Value = 1622.2222
$ws.Range("I82").Value = 1134
$ws.Range("J82").Value = 1866.3334
$ws.Range("K82").Value = 1134
$ws.Range("L82").Value = 1866.3334
$ws.Range("M82").Value = -773
$ws.Range("N82").Value = -2588.3334
$ws.Range("H85").Value = 1622.2222
$ws.Range("I85").Value = 1134
$ws.Range("J85").Value = 1866.3334
$ws.Range("K85").Value = 1134
$ws.Range("L85").Value = 1866.3334
$ws.Range("M85").Value = 114
$ws.Range("N85").Value = -4362.3334
$ws.Range("H93").Value = 1124.1177
$ws.Range("I93").Value = 960.75
$ws.Range("K93").Value = 960.75
$ws.Range("M93").Value = 287.25
$ws.Range("H100").Value = 4921.3335
$ws.Range("I100").Value = 4921.3335
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4921.3335
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -4380.3335
$ws.Range("H136").Value = 8829
$ws.Range("I136").Value = 8799
$ws.Range("K136").Value = 26397
$ws.Range("M136").Value = -23847

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2073.25
$ws.Range("I122").Value = 2079.9092
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6239.7276
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3789.7276
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 4200.5454
$ws.Range("I132").Value = 3030.2856
$ws.Range("J132").Value = 6248.5
$ws.Range("K132").Value = 9090.856800000001
$ws.Range("L132").Value = 18745.5
$ws.Range("M132").Value = -6560.856800000001
$ws.Range("N132").Value = -23805.5
$ws.Range("H136").Value = 1488.6471
$ws.Range("I136").Value = 976.8333
$ws.Range("K136").Value = 2930.4999
$ws.Range("M136").Value = -380.4998999999998
